# HRMS Regression test script First Commit
# Updates the Induction schedule: new In/Out time values, new employee
# name, highlights the In Time / Out Time / SearchTextBox headers, and
# moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data updates ---
# "In Time(U)" value: 10 -> 08
$ws.Range("D2").Value = "08"
# "Out Time(U)" value: 11 -> 09
$ws.Range("E2").Value = "09"
# "Employee Name" value: Gaikwad Ravina -> ritesh pandey
$ws.Range("G2").Value = "ritesh pandey"

# --- Formatting updates ---
# Highlight the In Time / Out Time / SearchTextBox header cells with a
# new fill color (indexed palette color 53).
$ws.Range("D1").Interior.ColorIndex = 53
$ws.Range("E1").Interior.ColorIndex = 53
$ws.Range("G1").Interior.ColorIndex = 53

# Touch the Employee Name value cell's font so its format is explicitly
# recorded (matches the author re-saving the cell format).
$ws.Range("G2").Font.Name = "Times New Roman"

# --- Selection ---
$ws.Range("E6").Select()
